$d = $word.ActiveDocument

# Locate the target paragraph (the one about the search algorithm).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*breadth first search algorithm*") {
        $target = $p
        break
    }
}

$paraStart = $target.Range.Start

$oldFirstRun = "Finally, using the breadth first search algorithm, I was able to "
$newFirstRun = "Finally, using the uniform cost search algorithm, I was able to "
$oldFirstRunLen = $oldFirstRun.Length
$newFirstRunLen = $newFirstRun.Length

# The original runs that make up the rest of the paragraph (after the
# sentence we are editing), captured here as plain text so their lengths -
# and therefore their boundaries - can be computed without relying on
# hard-coded absolute character offsets.
$tailRunTexts = @(
    "get the required output",
    " after a user enters ",
    "an input of a city and country",
    " for both start and finish destination",
    ". The output produced ",
    "consisted of",
    " an airline",
    " code, ",
    "an airport code ant the number of stops made",
    ". ",
    "It also consisted of the total number of flights needed to get to a destination as w",
    "ell as the number of stops that ",
    "have to",
    " be made. "
)

# Step 1: swap "breadth first" for "uniform cost" across the whole first
# run's range. This is a plain text replace, so the engine may coalesce this
# run with any immediately-following runs that share identical formatting
# (it does - the whole paragraph collapses into one run when this happens).
$r1End = $paraStart + $oldFirstRunLen
$r1 = $d.Range($paraStart, $r1End)
$r1.Text = $newFirstRun

# Step 2: re-introduce run boundaries. Toggling a (no-visible-effect-once-
# reverted) character formatting property on a sub-range forces the engine
# to materialize it as its own run without touching the text - so we can
# recreate every boundary that existed before step 1 (shifted by however
# much the edited sentence grew/shrank), plus the two new boundaries needed
# to carve "uniform cost" out on its own.
function Split-Run($start, $end) {
    $r = $d.Range($start, $end)
    $r.Bold = 1
    $r.Bold = 0
}

$seg1Len = "Finally, using the ".Length
$seg2Len = "uniform cost".Length
$p1 = $paraStart + $seg1Len
$p2 = $p1 + $seg2Len
$p3 = $paraStart + $newFirstRunLen

Split-Run $paraStart $p1
Split-Run $p1 $p2
Split-Run $p2 $p3

$cursor = $p3
foreach ($t in $tailRunTexts) {
    $tStart = $cursor
    $tEnd = $cursor + $t.Length
    Split-Run $tStart $tEnd
    $cursor = $tEnd
}

Write-Output "done"
